# Update the Option sheet's slot-type list:
# remove AMAZON.Actor / AdministrativeArea / AggregateRating / Airline,
# add AMAZON.EmailAddress, and refresh the slotType named range + active tab.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Option")

$slotTypes = @(
    "AMAZON.Airport","AMAZON.Animal","AMAZON.Artist","AMAZON.AT_CITY","AMAZON.AT_REGION","AMAZON.Athlete","AMAZON.Author","AMAZON.Book",
    "AMAZON.BookSeries","AMAZON.BroadcastChannel","AMAZON.CivicStructure","AMAZON.Color","AMAZON.Comic","AMAZON.Corporation","AMAZON.Country","AMAZON.CreativeWorkType",
    "AMAZON.DATE","AMAZON.DayOfWeek","AMAZON.DE_CITY","AMAZON.DE_FIRST_NAME","AMAZON.DE_REGION","AMAZON.Dessert","AMAZON.DeviceType","AMAZON.Director",
    "AMAZON.Drink","AMAZON.DURATION","AMAZON.EducationalOrganization","AMAZON.EmailAddress","AMAZON.EUROPE_CITY","AMAZON.EventType","AMAZON.Festival","AMAZON.FictionalCharacter",
    "AMAZON.FinancialService","AMAZON.Food","AMAZON.FoodEstablishment","AMAZON.FOUR_DIGIT_NUMBER","AMAZON.Game","AMAZON.GB_CITY","AMAZON.GB_FIRST_NAME","AMAZON.GB_REGION",
    "AMAZON.Genre","AMAZON.Landform","AMAZON.LandmarksOrHistoricalBuildings","AMAZON.Language","AMAZON.LocalBusiness","AMAZON.LocalBusinessType","AMAZON.MedicalOrganization","AMAZON.Month",
    "AMAZON.Movie","AMAZON.MovieSeries","AMAZON.MovieTheater","AMAZON.MusicAlbum","AMAZON.MusicCreativeWorkType","AMAZON.MusicEvent","AMAZON.MusicGroup","AMAZON.Musician",
    "AMAZON.MusicPlaylist","AMAZON.MusicRecording","AMAZON.MusicVenue","AMAZON.MusicVideo","AMAZON.NUMBER","AMAZON.Organization","AMAZON.Percentage","AMAZON.Person",
    "AMAZON.PhoneNumber","AMAZON.PostalAddress","AMAZON.Professional","AMAZON.ProfessionalType","AMAZON.RadioChannel","AMAZON.Residence","AMAZON.Room","AMAZON.ScreeningEvent",
    "AMAZON.Service","AMAZON.SocialMediaPlatform","AMAZON.SoftwareApplication","AMAZON.SoftwareGame","AMAZON.SpeedUnit","AMAZON.Sport","AMAZON.SportsEvent","AMAZON.SportsTeam",
    "AMAZON.StreetAddress","AMAZON.TelevisionChannel","AMAZON.TIME","AMAZON.TVEpisode","AMAZON.TVSeason","AMAZON.TVSeries","AMAZON.US_CITY","AMAZON.US_FIRST_NAME",
    "AMAZON.US_LAST_NAME","AMAZON.US_STATE","AMAZON.VideoGame","AMAZON.WeatherCondition","AMAZON.WeightUnit","AMAZON.WrittenCreativeWorkType"
)

# Clear the old rows (B2:B98) before writing the refreshed, shorter list
$ws.Range("B2:B98").ClearContents()

for ($i = 0; $i -lt $slotTypes.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $slotTypes[$i]
}

# Update the slotType defined name to the new, shorter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "slotType") {
        $n.RefersTo = "=Option!`$B`$2:`$B`$" + (1 + $slotTypes.Length)
    }
}

# Make the Option sheet the active tab/selection (matches the saved view state)
$ws.Activate()
$ws.Range("B" + (1 + $slotTypes.Length)).Select()

Write-Output "Updated slotType list to $($slotTypes.Length) entries"
